# "Sweden 1div Norra" bases update (28-05-2024 07:50)
#
# Several fixture rows were duplicated with their home/away legs stored in
# the wrong physical row (row N holds what should be row N+1's data, and
# vice versa). The fix swaps the full record - match id, teams, score
# components and every odds column - between each pair of rows, leaving
# the row-scoped columns (A = running id, C = league name, D = kickoff
# date) untouched since those stay tied to the physical row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($r1, $r2) {
    # Columns B and E..AD carry the per-match data that needs to trade
    # places; A, C and D stay put (id counter / league / date).
    $cols = @(
        'B',
        'E','F','G','H','I','J','K',
        'L','M','N','O','P','Q','R','S','T','U','V','W','X','Y','Z',
        'AA','AB','AC','AD'
    )

    foreach ($col in $cols) {
        $cell1 = $ws.Range("$col$r1")
        $cell2 = $ws.Range("$col$r2")
        $v1 = $cell1.Value()
        $v2 = $cell2.Value()
        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}

Swap-RowData 26 27
Swap-RowData 86 87
Swap-RowData 106 107
Swap-RowData 175 176
Swap-RowData 210 212
Swap-RowData 227 228
